$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 17, shifting the existing row 17
# ("HOLY FAMILY MARONITE CHURCH") down to row 18.
$ws.Rows.Item(17).Insert()
$ws.Rows.Item(17).RowHeight = $ws.Rows.Item(16).RowHeight

$ws.Range("A17").Value = "BEHRMAV ENTERPRISES LLC"
$ws.Range("B17").Value = "Larsen, Rick J"
$ws.Range("C17").Value = "015"
$ws.Range("E17").Value = "0008337"
